# Loan RBI, Variable Instalments
#
# This script reproduces the author's edit:
#  - A new (blank) column is inserted into the "Repayment schedule" sheet,
#    immediately before the old column N ("Late"). This pushes the old
#    N/O/P columns ("Late" / "heading" / "Outstanding") one slot to the
#    right (N->O, O->P, P->Q), and the freshly inserted column N inherits
#    the column width from column M (standard Excel "insert column"
#    behaviour of copying format from the column to the left).
#  - The "Repayment schedule" sheet becomes the active/selected sheet
#    (previously it was "Transactions"), with cell R8 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting the
# existing N/O/P columns right by one.
$ws.Columns("N").Insert()

# The newly inserted column takes on the width of column M, matching
# Excel's default "format same as left" behaviour on column insert.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with R8 selected
# (previously "Transactions" was active/selected).
$ws.Activate()
$ws.Range("R8").Select()
